$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A2 used to hold the shared string "Text2" - replace it with a plain number.
$ws.Range("A2").Value = 12345

# New row: a number cell styled with a brand-new "Century" font (this is
# what introduces the 3rd font / 3rd cellXfs entry in styles.xml).
$ws.Range("A3").Value = 123.45
$ws.Range("A3").Font.Name = "Century"
